$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Okada"
$ws.Range("C3").Value = 13494
$ws.Range("C11").Value = 9800

$ws.Range("B2").Select()
